$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed "Order History" export - new orders placed, table re-sorted with
# most recent Bill No first (rows 2-9).
# Columns: A=Bill No, B=Date, C=Cashier, D=KOT, E=Price, F=SGST, G=CGST, H=Tax, I=Food Items
$data = @(
    @(13, 45693.22928240741, "Reevan", 121, 50, 1, 1, 2, "Mango Lassi (x9)"),
    @(12, 45693.22928240741, "Karthik", 12, 450, 8, 8, 16, "Vanilla Shake (x4)"),
    @(11, 45693.22928240741, "Reevan", 111, 50, 1, 1, 2, "Veg Cheese Pops (x1)"),
    @(10, 45693.22928240741, "Karthik", 11, 450, 8, 8, 16, "Strawberry Shake (x3)"),
    @(9, 45692.22928240741, "Ajay Francis Anchan", 10, 870, 0, 0, 0, "Chicken Wrap (x9)"),
    @(8, 45692.22928240741, "Notsla Daniel", 9, 800, 5, 36, 11, "Butterscotch Lassi (x8)"),
    @(7, 45692.22928240741, "Ajay Francis Anchan", 8, 900, 5, 6, 11, "Chicken Burger (x3), Butterscotch Lassi (x7)"),
    @(1, 45687.22928240741, "Ajay Francis Anchan", 1, 130, 2, 3, 5, "Chicken Burger (x3), Mango Lassi (x1)")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $row++
}
